$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').NumberFormat = '@'
$ws.Range('D2').Value = '27.239.84'
$ws.Range('D2').Style = 'Normal'
$ws.Range('E2').Value = '  +0.30%  '
$ws.Range('D3').NumberFormat = '@'
$ws.Range('D3').Value = '1.906.17'
$ws.Range('D3').Style = 'Normal'
$ws.Range('D4').NumberFormat = '@'
$ws.Range('D4').Value = '1.001'
$ws.Range('D4').Style = 'Normal'
$ws.Range('E4').Value = '  +0.04%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '307.37'
$ws.Range('D5').Style = 'Normal'
$ws.Range('E5').Value = '  +0.24%  '
$ws.Range('E6').Value = '  +0.12%  '
$ws.Range('D7').NumberFormat = '@'
$ws.Range('D7').Value = '0.5255'
$ws.Range('D7').Style = 'Normal'
$ws.Range('E7').Value = '  +0.38%  '
$ws.Range('D8').NumberFormat = '@'
$ws.Range('D8').Value = '0.3813'
$ws.Range('D8').Style = 'Normal'
$ws.Range('E8').Value = '  +1.34%  '
$ws.Range('D9').NumberFormat = '@'
$ws.Range('D9').Value = '0.07293'
$ws.Range('D9').Style = 'Normal'
$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '21.74'
$ws.Range('D10').Style = 'Normal'
$ws.Range('E10').Value = '  +2.59%  '
$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '0.9037'
$ws.Range('D11').Style = 'Normal'
$ws.Range('E11').Value = '  -0.34%  '
$ws.Range('D12').NumberFormat = '@'
$ws.Range('D12').Value = '0.08207'
$ws.Range('D12').Style = 'Normal'
$ws.Range('E12').Value = '  -3.27%  '
$ws.Range('D13').NumberFormat = '@'
$ws.Range('D13').Value = '96.35'
$ws.Range('D13').Style = 'Normal'
$ws.Range('E13').Value = '  -0.65%  '
$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '5.357'
$ws.Range('D14').Style = 'Normal'
$ws.Range('E14').Value = '  +1.18%  '
$ws.Range('D15').NumberFormat = '@'
$ws.Range('D15').Value = '1.500.26'
$ws.Range('D15').Style = 'Normal'
$ws.Range('E15').Value = '  -21.88%  '
$ws.Range('D16').NumberFormat = '@'
$ws.Range('D16').Value = '1.000'
$ws.Range('D16').Style = 'Normal'
$ws.Range('E16').Value = '  -0.07%  '
$ws.Range('D17').NumberFormat = '@'
$ws.Range('D17').Value = '0.000008659'
$ws.Range('D17').Style = 'Normal'
$ws.Range('E17').Value = '  -0.47%  '
$ws.Range('D18').NumberFormat = '@'
$ws.Range('D18').Value = '14.77'
$ws.Range('D18').Style = 'Normal'
$ws.Range('E18').Value = '  +1.50%  '
$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '27.269.42'
$ws.Range('D20').Style = 'Normal'
$ws.Range('E20').Value = '  +0.28%  '
$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '5.111'
$ws.Range('D21').Style = 'Normal'
$ws.Range('E21').Value = '  +0.49%  '
$ws.Range('E22').Value = '  +1.70%  '
$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '6.505'
$ws.Range('D23').Style = 'Normal'
$ws.Range('E23').Value = '  +1.03%  '
$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '2.343'
$ws.Range('D24').Style = 'Normal'
$ws.Range('E24').Value = '  +1.07%  '
$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '150.28'
$ws.Range('D25').Style = 'Normal'
$ws.Range('E25').Value = '  +2.31%  '
$ws.Range('E26').Value = '  +0.00%  '
$ws.Range('E27').Value = '  -0.26%  '
$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '116.66'
$ws.Range('D28').Style = 'Normal'
$ws.Range('E28').Value = '  +1.29%  '
$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '4.845'
$ws.Range('D29').Style = 'Normal'
$ws.Range('E29').Value = '  +0.34%  '
$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '4.858'
$ws.Range('D30').Style = 'Normal'
$ws.Range('E30').Value = '  -1.29%  '
$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '0.09252'
$ws.Range('D31').Style = 'Normal'
$ws.Range('E31').Value = '  -0.59%  '
$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '0.8313'
$ws.Range('D32').Style = 'Normal'
$ws.Range('E32').Value = '  +3.64%  '
$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '0.05054'
$ws.Range('D33').Style = 'Normal'
$ws.Range('E33').Value = '  -0.10%  '
$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '1.229'
$ws.Range('D34').Style = 'Normal'
$ws.Range('E34').Value = '  -1.31%  '
$ws.Range('E35').Value = '  +1.36%  '
$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '3.354'
$ws.Range('D36').Style = 'Normal'
$ws.Range('E36').Value = '  -2.49%  '
$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '2.735'
$ws.Range('D37').Style = 'Normal'
$ws.Range('E37').Value = '  +5.17%  '
$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '0.5782'
$ws.Range('D38').Style = 'Normal'
$ws.Range('E38').Value = '  +1.19%  '
$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '0.02005'
$ws.Range('D39').Style = 'Normal'
$ws.Range('E39').Value = '  +0.12%  '
$ws.Range('E40').Value = '  +0.48%  '
$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '9.226'
$ws.Range('D41').Style = 'Normal'
$ws.Range('E41').Value = '  +1.11%  '
$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '6.605'
$ws.Range('D42').Style = 'Normal'
$ws.Range('E42').Value = '  -0.37%  '
$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '115.81'
$ws.Range('D43').Style = 'Normal'
$ws.Range('E43').Value = '  -0.14%  '
$ws.Range('E44').Value = '  +0.37%  '
$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '0.4918'
$ws.Range('D45').Style = 'Normal'
$ws.Range('E45').Value = '  +1.12%  '
$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '10.24'
$ws.Range('D46').Style = 'Normal'
$ws.Range('E46').Value = '  +0.88%  '
$ws.Range('E47').Value = '  +0.13%  '
$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '1.644'
$ws.Range('D48').Style = 'Normal'
$ws.Range('E48').Value = '  +1.23%  '
$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '38.94'
$ws.Range('D49').Style = 'Normal'
$ws.Range('E49').Value = '  +3.07%  '
$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '64.41'
$ws.Range('D50').Style = 'Normal'
$ws.Range('E50').Value = '  +0.28%  '
$ws.Range('E51').Value = '  +1.63%  '
